# R33 was found to cause a gain error in the transimpedance amplifier of the
# first V0.2 prototypes. Move its designator from the 20k resistor group to
# the 12k resistor group in the BOM.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep gridlines shown (avoids the sheet view losing its "show gridlines"
# state when the view block gets regenerated below).
$excel.ActiveWindow.DisplayGridlines = $true

# 12k Resistor, 1% group (row 13) gains R33
$ws.Range("B13").Value = "R1,R2,R26,R28,R34,R36,R33"

# 20k Resistor, 1% group (row 14) loses R33
$ws.Range("B14").Value = "R3,R5,R6,R8,R10,R13,R23,R25,R29"

# Widening column B was needed to fit the longer designator text; re-assert
# the other custom column widths so they keep their explicit (customWidth)
# state instead of reverting to sheet defaults.
$ws.Columns.Item(1).ColumnWidth = 65.857142857142857
$ws.Columns.Item(2).ColumnWidth = 50.285714285714285
$ws.Columns.Item(3).ColumnWidth = 34.0
$ws.Columns.Item(4).ColumnWidth = 27.285714285714285
$ws.Columns.Item(5).ColumnWidth = 11.0

# Move the saved selection/active cell to B14 (the edited cell)
[void]$ws.Range("B14").Select()
